# MINERVA 시스템 구성도 - slide 2 updates
#
# EMU -> point conversion: 1 pt = 12700 EMU. PowerPoint's Shape.Left/Top/
# Width/Height COM properties are expressed in points, so every EMU offset
# from the target XML is divided by 12700.0 (kept as a full-precision
# float expression) to round-trip back to the exact EMU value on save.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1) "카메라 (CCTV)" shape - merge the split "(CCTV" / ")" runs into one
# ---------------------------------------------------------------------
$cctvShape = $s.Shapes.Item(1)
$cctvPara = $cctvShape.TextFrame.TextRange.Paragraphs(2)
$cctvPara.Text = "temp"
$cctvPara.Text = "(CCTV)"

# ---------------------------------------------------------------------
# 2) "실시간 객체 인식 (YOLOv8 + Python)" shape - reposition + merge runs
# ---------------------------------------------------------------------
$yoloShape = $s.Shapes.Item(5)
$yoloShape.Left = 3848143 / 12700.0
$yoloShape.Top = 3997202 / 12700.0

$yoloPara = $yoloShape.TextFrame.TextRange.Paragraphs(2)
$yoloPara.Text = "temp"
$yoloPara.Text = "(YOLOv8 + Python)"

# ---------------------------------------------------------------------
# 3) "관리자" shape - reposition & resize
# ---------------------------------------------------------------------
$managerShape = $s.Shapes.Item(8)
$managerShape.Left = 669716 / 12700.0
$managerShape.Top = 2626854 / 12700.0
$managerShape.Width = 2046408 / 12700.0
$managerShape.Height = 1174200 / 12700.0

# ---------------------------------------------------------------------
# 4) "실시간 객체 인식" shape (big box, top) - reposition + merge
#    "현황 " / "분석" runs into a single "현황 분석" run
# ---------------------------------------------------------------------
$realtimeShape = $s.Shapes.Item(10)
$realtimeShape.Left = 3680212 / 12700.0
$realtimeShape.Top = 1616497 / 12700.0

$statusPara = $realtimeShape.TextFrame.TextRange.Paragraphs(5)
$statusPara.Text = "temp"
$statusPara.Text = "현황 분석"

# ---------------------------------------------------------------------
# 5) Three new label shapes, duplicated from the existing plain
#    rectangle "직사각형 10" (id 11) so they inherit the same
#    rect/noFill-line/p:style formatting.
# ---------------------------------------------------------------------
$template = $s.Shapes.Item(4)

$webShape = $template.Duplicate().Item(1)
$webShape.Name = "직사각형 11"
$webShape.Left = 5325952 / 12700.0
$webShape.Top = 1376290 / 12700.0
$webShape.Width = 1117880 / 12700.0
$webShape.Height = 320808 / 12700.0
$webShape.TextFrame.TextRange.Text = "웹"

$aiModuleShape = $template.Duplicate().Item(1)
$aiModuleShape.Name = "직사각형 14"
$aiModuleShape.Left = 5325952 / 12700.0
$aiModuleShape.Top = 3801053 / 12700.0
$aiModuleShape.Width = 1117880 / 12700.0
$aiModuleShape.Height = 320808 / 12700.0
$aiModuleShape.TextFrame.TextRange.Text = "AI 모듈"

$userShape = $template.Duplicate().Item(1)
$userShape.Name = "직사각형 16"
$userShape.Left = 1133980 / 12700.0
$userShape.Top = 2466450 / 12700.0
$userShape.Width = 1117880 / 12700.0
$userShape.Height = 320808 / 12700.0
$userShape.TextFrame.TextRange.Text = "사용자"

Write-Host "done"
